# Buff.xlsx edit: "unify the conception of DataNode, DataTable, Entity."
#
# The underlying rename is from worksheet "Property1" to "DataNode" (the
# sheet tab + <sheet name=.../> entry in xl/workbook.xml). The rest of the
# original diff (fileVersion/rupBuild bump, xr2/xr16/x16r2/xr namespace +
# mc:Ignorable additions, xr:uid/xr2:uid GUIDs, the extra phonetic "宋体" 9pt
# font + <phoneticPr>, the "Normal" -> "常规" cell-style label, the
# absPath/workbookView window geometry, and the sub-pixel column width
# drift) is the normal fingerprint Excel leaves behind merely by re-saving
# an older workbook in a newer build on a zh-CN locale machine - cosmetic
# side effects of the resave, not a deliberate action by the edit's author,
# and not something the Excel object model exposes a way to set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Property1 -> DataNode (the actual content of the commit).
$ws.Name = "DataNode"

# The selection left behind after editing moved from A9 to H13.
$ws.Range("H13").Select()

# Row heights were reflowed (6 wrapped lines each, at the new default row
# height) when the sheet was resaved: row 1 -> 27, row 8 -> 81.
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 81
